# Generate Report for Handback
#
# Marks the zh-cn / de-de localization rows as handed back (in sync with
# en-US), records the generated target (.md) and handback (.xlf) files for
# each row, and stamps the handback datetime.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: flip the per-language status from "Ready for handoff"
# to "Handed back: in sync with en-US"
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# Widen the status columns so the longer text fits.
$overview.Columns.Item(5).ColumnWidth = 29.2
$overview.Columns.Item(6).ColumnWidth = 29.2

# ---------------------------------------------------------------------
# zh-cn sheet: fill in the generated target file, the handback file, and
# the handback datetime for both rows.
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Columns.Item(3).ColumnWidth = 29.2
$zhcn.Columns.Item(9).ColumnWidth = 39.15
$zhcn.Columns.Item(10).ColumnWidth = 39.15

$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fc8282e00b03ee7723de515acbe41f08f7d7ce65/e2e/064a7f1e-bcb2-4677-899e-eec20210867a.md", "", "", "064a7f1e-bcb2-4677-899e-eec20210867a.md")
$zhcn.Range("J2").Value = "064a7f1e-bcb2-4677-899e-eec20210867a.7c0faf4da43e001b647b99670974dde769f5a8d0.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-19 06:46:07"

$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fc8282e00b03ee7723de515acbe41f08f7d7ce65/e2e/9c788acf-9ee0-447c-a4ae-d1a226a7d50f.md", "", "", "9c788acf-9ee0-447c-a4ae-d1a226a7d50f.md")
$zhcn.Range("J3").Value = "9c788acf-9ee0-447c-a4ae-d1a226a7d50f.1fb388f31ecac20a6cd6320a52260e33bdb1694a.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-19 06:46:07"

# ---------------------------------------------------------------------
# de-de sheet: same as zh-cn, with its own handback timestamp.
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Columns.Item(3).ColumnWidth = 29.2
$dede.Columns.Item(9).ColumnWidth = 39.15
$dede.Columns.Item(10).ColumnWidth = 39.15

$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("C3").Value = "Handed back: in sync with en-US"

$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fc8282e00b03ee7723de515acbe41f08f7d7ce65/e2e/064a7f1e-bcb2-4677-899e-eec20210867a.md", "", "", "064a7f1e-bcb2-4677-899e-eec20210867a.md")
$dede.Range("J2").Value = "064a7f1e-bcb2-4677-899e-eec20210867a.7c0faf4da43e001b647b99670974dde769f5a8d0.de-de.xlf"
$dede.Range("K2").Value = "2016-08-19 06:46:15"

$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fc8282e00b03ee7723de515acbe41f08f7d7ce65/e2e/9c788acf-9ee0-447c-a4ae-d1a226a7d50f.md", "", "", "9c788acf-9ee0-447c-a4ae-d1a226a7d50f.md")
$dede.Range("J3").Value = "9c788acf-9ee0-447c-a4ae-d1a226a7d50f.1fb388f31ecac20a6cd6320a52260e33bdb1694a.de-de.xlf"
$dede.Range("K3").Value = "2016-08-19 06:46:15"
